$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Maggie  Farrell -> Chrissy Cummings
$ws.Range("A3").Value = "Chrissy Cummings"
$ws.Range("B3").Value = 29
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = 6

# Row 4: Makayla Baca -> Danielle Mai
$ws.Range("A4").Value = "Danielle Mai"
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 1

# Row 5: Justyne Martinez -> Jasmine Saiz
$ws.Range("A5").Value = "Jasmine Saiz"
$ws.Range("B5").Value = 56
$ws.Range("C5").Value = 27
$ws.Range("D5").Value = 10

# Row 6: Vy Torino -> Karen Trevizo
$ws.Range("A6").Value = "Karen Trevizo"
$ws.Range("B6").Value = 27
$ws.Range("C6").Value = 21
$ws.Range("D6").Value = 1
